$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 145; existing rows 145-174 shift down to 146-175.
$ws.Rows.Item(145).Insert()

# Populate the newly inserted row 145 with the new record's data.
$ws.Range("A145").Value = 5
$ws.Range("B145").Value = "Macroferia Regional de Talca"
$ws.Range("C145").Value = "Maule"
$ws.Range("D145").Value = 45275
$ws.Range("E145").Value = 7
$ws.Range("F145").Value = 100112022
$ws.Range("G145").Value = "Arveja Verde"
$ws.Range("H145").Value = "Sin especificar"
$ws.Range("I145").Value = "Primera"
$ws.Range("J145").Value = 300
$ws.Range("K145").Value = 23000
$ws.Range("L145").Value = 23000
$ws.Range("M145").Value = 23000
$ws.Range("N145").Value = "$/saco 25 kilos"
$ws.Range("O145").Value = "Región del Maule"
$ws.Range("P145").Value = 920
$ws.Range("Q145").Value = 25
$ws.Range("R145").Value = "Hortaliza"
